$wb = $excel.ActiveWorkbook

# --- TestSuite sheet: flip OpenAccount run mode from "n" to "y" ---
$ws1 = $wb.Worksheets.Item("TestSuite")
$ws1.Range("B4").Value = "y"

# --- OpenAccount sheet: add a RunMode column (D) driving each row ---
$ws3 = $wb.Worksheets.Item("OpenAccount")
$ws3.Range("D1").Value = "RunMode"
$ws3.Range("D2").Value = "y"
$ws3.Range("D3").Value = "y"
$ws3.Range("D4").Value = "y"
$ws3.Range("D5").Value = "n"

# --- Update the saved selections to match the edited cells ---
$ws1.Activate()
$ws1.Range("B4").Select()

$ws3.Activate()
$ws3.Range("D6").Select()

$ws1.Activate()
